$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.539
$ws.Range("C21").Value = -12.779
$ws.Range("C23").Value = -12.937
$ws.Range("C25").Value = -12.763
